# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# Rows 2-5: F2 6602->6615, F3 41->42, F4 191->193, F5 1027->1029

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 6615
    $ws.Range("F3").Value = 42
    $ws.Range("F4").Value = 193
    $ws.Range("F5").Value = 1029
}
